$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 396; $r++) {
    $ws.Cells.Item($r, 15).Value = "2023-01-15 12:54:58"
}
